$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '31.227.14'
Set-TextValue 'E2' '  +2.94%  '
Set-TextValue 'D3' '1.981.02'
Set-TextValue 'E3' '  +5.98%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '0.7931'
Set-TextValue 'E5' '  +68.63%  '
Set-TextValue 'D6' '252.85'
Set-TextValue 'E6' '  +3.79%  '
Set-TextValue 'D7' '1.001'
Set-TextValue 'E7' '  +0.07%  '
Set-TextValue 'D8' '0.3371'
Set-TextValue 'E8' '  +17.50%  '
Set-TextValue 'D9' '25.62'
Set-TextValue 'E9' '  +16.28%  '
Set-TextValue 'D10' '0.06916'
Set-TextValue 'E10' '  +7.28%  '
Set-TextValue 'D11' '0.8323'
Set-TextValue 'E11' '  +15.62%  '
Set-TextValue 'D12' '0.08098'
Set-TextValue 'E12' '  +4.26%  '
Set-TextValue 'D13' '1.989.93'
Set-TextValue 'E13' '  +6.38%  '
Set-TextValue 'D14' '99.99'
Set-TextValue 'E14' '  +4.22%  '
Set-TextValue 'D15' '5.438'
Set-TextValue 'E15' '  +6.24%  '
Set-TextValue 'D16' '272.97'
Set-TextValue 'E16' '  -1.86%  '
Set-TextValue 'D17' '31.245.54'
Set-TextValue 'E17' '  +3.05%  '
Set-TextValue 'D18' '13.81'
Set-TextValue 'E18' '  +6.64%  '
Set-TextValue 'D19' '0.000007906'
Set-TextValue 'E19' '  +5.40%  '
Set-TextValue 'D20' '2.249.55'
Set-TextValue 'E20' '  +6.35%  '
Set-TextValue 'D21' '5.705'
Set-TextValue 'E21' '  +9.43%  '
Set-TextValue 'E22' '  +0.25%  '
Set-TextValue 'D23' '1.001'
Set-TextValue 'E23' '  +0.07%  '
Set-TextValue 'E24' '  +11.32%  '
Set-TextValue 'D25' '9.610'
Set-TextValue 'E25' '  +6.44%  '
Set-TextValue 'D26' '164.49'
Set-TextValue 'E26' '  +0.74%  '
Set-TextValue 'D27' '0.1470'
Set-TextValue 'E27' '  +53.32%  '
Set-TextValue 'D28' '19.74'
Set-TextValue 'E28' '  +5.93%  '
Set-TextValue 'D29' '2.176'
Set-TextValue 'E29' '  +16.36%  '
Set-TextValue 'B30' 'PancakeSwap'
Set-TextValue 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '1.564'
Set-TextValue 'E30' '  +6.23%  '
Set-TextValue 'B31' 'Toncoin'
Set-TextValue 'C31' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '1.356'
Set-TextValue 'E31' '  +2.74%  '
Set-TextValue 'D32' '4.547'
Set-TextValue 'E32' '  +8.34%  '
Set-TextValue 'E33' '  +5.57%  '
Set-TextValue 'D34' '0.05142'
Set-TextValue 'E34' '  +7.20%  '
Set-TextValue 'D35' '1.205'
Set-TextValue 'E35' '  +7.90%  '
Set-TextValue 'D36' '0.7519'
Set-TextValue 'E36' '  +9.29%  '
Set-TextValue 'D37' '2.789'
Set-TextValue 'E37' '  +2.87%  '
Set-TextValue 'B38' 'VeChain'
Set-TextValue 'C38' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D38' '0.02002'
Set-TextValue 'E38' '  +6.86%  '
Set-TextValue 'B39' 'MXToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D39' '2.909'
Set-TextValue 'E39' '  +3.63%  '
Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '6.595'
Set-TextValue 'E40' '  +6.38%  '
Set-TextValue 'B41' 'Aave'
Set-TextValue 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '78.05'
Set-TextValue 'E41' '  +5.31%  '
Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.4632'
Set-TextValue 'E42' '  +10.09%  '
Set-TextValue 'B43' 'RenderToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D43' '2.049'
Set-TextValue 'E43' '  +6.15%  '
Set-TextValue 'B44' 'TrustWalletToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D44' '0.8506'
Set-TextValue 'E44' '  +2.72%  '
Set-TextValue 'B45' 'Quant'
Set-TextValue 'C45' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D45' '105.17'
Set-TextValue 'E45' '  +4.68%  '
Set-TextValue 'B46' 'PaxDollar'
Set-TextValue 'C46' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D46' '1.001'
Set-TextValue 'E46' '  +0.19%  '
Set-TextValue 'B47' 'EnergySwap'
Set-TextValue 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '9.967'
Set-TextValue 'E47' '  +4.44%  '
Set-TextValue 'B48' 'Aptos'
Set-TextValue 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D48' '7.463'
Set-TextValue 'E48' '  +7.54%  '
Set-TextValue 'B49' 'Elrond'
Set-TextValue 'C49' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D49' '36.35'
Set-TextValue 'E49' '  +3.25%  '
Set-TextValue 'B50' 'Decentraland'
Set-TextValue 'C50' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D50' '0.4258'
Set-TextValue 'E50' '  +8.93%  '
Set-TextValue 'B51' 'Maker'
Set-TextValue 'C51' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D51' '921.48'
Set-TextValue 'E51' '  +2.38%  '
